$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P0001")
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 75
Write-Host "zoom" $win.Zoom
